$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 4")

$ws.Range("D3").Value  = "0.34 [0.24 - 0.45]"
$ws.Range("D4").Value  = "0.78 [0.66 - 0.89]"
$ws.Range("D5").Value  = "0.68 [0.61 - 0.75]"
$ws.Range("D7").Value  = "0.37 [0.27 - 0.47]"
$ws.Range("D11").Value = "0.38 [0.27 - 0.47]"
$ws.Range("D13").Value = "0.69 [0.61 - 0.75]"
